# Update "想去人数" (want-to-go count) values on the "展览" and "全部类型" sheets
# to match the newly scraped data (output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 286
$ws1.Range("F4").Value = 2329
$ws1.Range("F5").Value = 1745
$ws1.Range("F6").Value = 334
$ws1.Range("F7").Value = 96
$ws1.Range("F8").Value = 804
$ws1.Range("F9").Value = 164

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 286
$ws4.Range("F4").Value = 2329
$ws4.Range("F5").Value = 1745
$ws4.Range("F6").Value = 334
$ws4.Range("F8").Value = 96
$ws4.Range("F9").Value = 804
$ws4.Range("F10").Value = 164
